# Update gh-pages to output generated at 456a3b4
# Increment "想去人数" (F column) counts on the 展览, 演出, and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 14245
$wsExhibit.Range("F6").Value = 555
$wsExhibit.Range("F7").Value = 1489
$wsExhibit.Range("F8").Value = 142

# 演出 (Performances) sheet
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 11

# 全部类型 (All types) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14245
$wsAll.Range("F6").Value = 11
$wsAll.Range("F8").Value = 555
$wsAll.Range("F9").Value = 1489
$wsAll.Range("F11").Value = 142

$wb.Save()
